$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 0
$ws.Range("F6").Value = 5
$ws.Range("F9").Value = -3
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = 6
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = -5
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = -5
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = -1
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = -1
$ws.Range("F27").Value = -2
$ws.Range("F29").Value = -2
$ws.Range("F32").Value = -11
$ws.Range("F34").Value = -3
$ws.Range("F35").Value = 0
$ws.Range("F37").Value = -3
$ws.Range("F42").Value = -2
$ws.Range("F43").Value = -3
$ws.Range("F44").Value = -5
$ws.Range("F46").Value = -4
$ws.Range("F47").Value = -6
$ws.Range("F48").Value = 2
$ws.Range("F49").Value = 2
$ws.Range("F50").Value = -1
$ws.Range("F52").Value = -5
